$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Update the letter date: "September 19, 2025" -> "September 21, 2025".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Split the first mailing-address paragraph
#    "919 Story Road, San Jose CA 95122" into two separate paragraphs:
#    "919 Story Road" and "San Jose, CA 95122".
# ---------------------------------------------------------------------------
$streetLine = "919 Story Road"
$cityLine = "San Jose, CA 95122"
$fullAddress = "$streetLine, San Jose CA 95122"

$addrIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "$fullAddress*") {
        $addrIndex = $i
        break
    }
}

if ($addrIndex -gt 0) {
    $addrRange = $d.Paragraphs.Item($addrIndex).Range
    $splitPoint = $addrRange.Start + $streetLine.Length
    $splitRange = $d.Range($splitPoint, $splitPoint)
    $splitRange.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($addrIndex + 1)
    $newPara.Range.Text = $cityLine
}

# ---------------------------------------------------------------------------
# 3. Remove the now-redundant blank "No Spacing" paragraph that follows the
#    "...Board of Directors" signature line. The paragraph index is located
#    dynamically since the insert above shifts numbering for everything
#    that follows it.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Board of Directors*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0 -and $targetIndex -lt $count) {
    $blankPara = $d.Paragraphs.Item($targetIndex + 1)
    $blankText = $blankPara.Range.Text.Trim([char]13, [char]7, ' ')
    if ($blankText -eq "" -and $blankPara.Style.NameLocal -eq "No Spacing") {
        $blankPara.Range.Delete()
    }
}
